$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generation")

# Updated "RES" (column G) generation figures for the final database run + TEMBA tests
$ws.Range("G2").Value  = 367.15391229578677
$ws.Range("G3").Value  = 18.056749785038694
$ws.Range("G4").Value  = 732.58813413585551
$ws.Range("G5").Value  = 12.037833190025795
$ws.Range("G6").Value  = 91.143594153052447
$ws.Range("G7").Value  = 379.19174548581253
$ws.Range("G9").Value  = 27.515047291487534
$ws.Range("G10").Value = 1177.9879621668099
$ws.Range("G12").Value = 773.86070507308682
$ws.Range("G13").Value = 68.787618228718827
$ws.Range("G14").Value = 0.85984522785898543
$ws.Range("G15").Value = 295.78675838349096
$ws.Range("G17").Value = 216.68099742046431
$ws.Range("G18").Value = 10.318142734307825
$ws.Range("G19").Value = 15.477214101461737
$ws.Range("G20").Value = 542.5623387790198
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 2.1496130696474633
$ws.Range("G24").Value = 4.2992261392949267
$ws.Range("G25").Value = 220.12037833190027
$ws.Range("G26").Value = 127.25709372312984

# Selection moved from a single cell to the full data range A1:B26
$ws.Range("A1:B26").Select()
